$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update salt (Magnesium chloride) row 8: set new static lower/midpoint/upper values
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Update salt (Zinc sulfate) row 9: set new static lower/midpoint/upper values
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# Reflect the selection left after the edit (rows 8:9 selected, active cell A8)
$ws.Rows("8:9").Select()
